$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift values C1/D1/E1 -> C1=prediction, D1=rejection-f, E1=max
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# --- Data rows (2..156): C becomes the family label (text), D stays the family label,
# --- E becomes numeric 1
$lastRow = 156
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Range("D$r").Value2
    $ws.Range("C$r").Value = $label
    $ws.Range("E$r").Value = 1
}
